# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorganiza la tabla de estado de cuenta (filas 16-38) de "agrupado por
# trabajador" a "agrupado por periodo de mora", y corrige el Salario Basico
# de NORELIS REYES BELTRAN (antes 1.423.500, ahora 877.803 para todos sus
# periodos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cada fila: Tipo Doc, N Doc Trabajador, Nombre Trabajador, Periodo Mora, Valor Mora, Salario Basico
$data = @(
    @("CC", "1045229424", "YARELIS VEGA JIMENEZ",          "1912", 33125, 877803),
    @("CC", "50988101",   "NORELIS REYES BELTRAN",         "1912", 33125, 877803),
    @("CC", "1045229424", "YARELIS VEGA JIMENEZ",          "2001", 35112, 877803),
    @("CC", "50988101",   "NORELIS REYES BELTRAN",         "2001", 35112, 877803),
    @("CC", "45493570",   "SANDRA BIBIANA VEGA HERNANDEZ", "2001", 35112, 877803),
    @("CC", "1045229424", "YARELIS VEGA JIMENEZ",          "2002", 35112, 877803),
    @("CC", "50988101",   "NORELIS REYES BELTRAN",         "2002", 35112, 877803),
    @("CC", "45493570",   "SANDRA BIBIANA VEGA HERNANDEZ", "2002", 35112, 877803),
    @("CC", "1045229424", "YARELIS VEGA JIMENEZ",          "2003", 35112, 877803),
    @("CC", "50988101",   "NORELIS REYES BELTRAN",         "2003", 35112, 877803),
    @("CC", "45493570",   "SANDRA BIBIANA VEGA HERNANDEZ", "2003", 35112, 877803),
    @("CC", "1045229424", "YARELIS VEGA JIMENEZ",          "2004", 35112, 877803),
    @("CC", "50988101",   "NORELIS REYES BELTRAN",         "2004", 35112, 877803),
    @("CC", "45493570",   "SANDRA BIBIANA VEGA HERNANDEZ", "2004", 35112, 877803),
    @("CC", "1045229424", "YARELIS VEGA JIMENEZ",          "2005", 35112, 877803),
    @("CC", "50988101",   "NORELIS REYES BELTRAN",         "2005", 35112, 877803),
    @("CC", "45493570",   "SANDRA BIBIANA VEGA HERNANDEZ", "2005", 35112, 877803),
    @("CC", "1045229424", "YARELIS VEGA JIMENEZ",          "2006", 35112, 877803),
    @("CC", "50988101",   "NORELIS REYES BELTRAN",         "2006", 35112, 877803),
    @("CC", "45493570",   "SANDRA BIBIANA VEGA HERNANDEZ", "2006", 35112, 877803),
    @("CC", "1045229424", "YARELIS VEGA JIMENEZ",          "2007", 25749, 877803),
    @("CC", "50988101",   "NORELIS REYES BELTRAN",         "2007", 25749, 877803),
    @("CC", "45493570",   "SANDRA BIBIANA VEGA HERNANDEZ", "2007", 25749, 877803)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
}
